$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Akurasi"

# Update confusion matrix values - first table (SVM Kernel Linear / SVM Kernel RBF)
$ws.Range("D4").Value = 260
$ws.Range("E4").Value = 68
$ws.Range("D5").Value = 28
$ws.Range("E5").Value = 1208

$ws.Range("I4").Value = 242
$ws.Range("J4").Value = 86
$ws.Range("I5").Value = 28
$ws.Range("J5").Value = 1208

# Update confusion matrix values - second table (KNN)
$ws.Range("D10").Value = 143
$ws.Range("E10").Value = 149
$ws.Range("D11").Value = 26
$ws.Range("E11").Value = 1246

# Update label for KNN table
$ws.Range("B7").Value = "KNN K = 4"

# Update accuracy summary values
$ws.Range("M2").Value = "0.938619 (93.9%)"
$ws.Range("M3").Value = "0.927109 (92.7%)"
$ws.Range("M4").Value = "0.888107 (88.8%)"

# Update selection
$ws.Range("J6").Select()
